$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 110.15625
$ws.Range("I33").Value = 101.96429
$ws.Range("J33").Value = 167.5
$ws.Range("K33").Value = 101.96429
$ws.Range("L33").Value = 167.5
$ws.Range("M33").Value = 127.03571
$ws.Range("N33").Value = -625.5
$ws.Range("H40").Value = 2000.2
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 2175.25
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 2175.25
$ws.Range("M40").Value = -1125
$ws.Range("N40").Value = -2525.25
$ws.Range("H127").Value = 2200.6765
$ws.Range("I127").Value = 705.1
$ws.Range("J127").Value = 2823.8333
$ws.Range("K127").Value = 2115.3
$ws.Range("L127").Value = 8471.499899999999
$ws.Range("M127").Value = 2844.7
$ws.Range("N127").Value = -18391.4999
$ws.Range("H135").Value = 726.92426
$ws.Range("I135").Value = 396.54718
$ws.Range("J135").Value = 2073.8462
$ws.Range("K135").Value = 3568.92462
$ws.Range("L135").Value = 18664.6158
$ws.Range("M135").Value = -1033.92462
$ws.Range("N135").Value = -23734.6158
$ws.Range("H137").Value = 730.1087
$ws.Range("I137").Value = 647.94446
$ws.Range("J137").Value = 782.9286
$ws.Range("K137").Value = 1943.83338
$ws.Range("L137").Value = 2348.7858
$ws.Range("M137").Value = 606.16662
$ws.Range("N137").Value = -7448.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 935.41
$ws.Range("I32").Value = 729.6322
$ws.Range("J32").Value = 2312.5386
$ws.Range("K32").Value = 729.6322
$ws.Range("L32").Value = 2312.5386
$ws.Range("M32").Value = -442.6322
$ws.Range("N32").Value = -2886.5386
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -586
$ws.Range("N41").ClearContents()
$ws.Range("H55").Value = 251134
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 251134
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 251134
$ws.Range("N55").Value = -251764
$ws.Range("H61").Value = 1014.25
$ws.Range("I61").Value = 847.8421
$ws.Range("J61").Value = 1365.5555
$ws.Range("K61").Value = 847.8421
$ws.Range("L61").Value = 1365.5555
$ws.Range("M61").Value = -635.8421
$ws.Range("N61").Value = -1789.5555
$ws.Range("H132").Value = 1584.409
$ws.Range("I132").Value = 1397.6666
$ws.Range("J132").Value = 2424.75
$ws.Range("K132").Value = 4192.9998
$ws.Range("L132").Value = 7274.25
$ws.Range("M132").Value = -1662.9998
$ws.Range("N132").Value = -12334.25
$ws.Range("H136").Value = 1014.25
$ws.Range("I136").Value = 847.8421
$ws.Range("J136").Value = 1365.5555
$ws.Range("K136").Value = 2543.5263
$ws.Range("L136").Value = 4096.666499999999
$ws.Range("M136").Value = 6.473700000000008
$ws.Range("N136").Value = -9196.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18645.17
$ws.Range("I134").Value = 1477.2291
$ws.Range("J134").Value = 93559.82
$ws.Range("K134").Value = 4431.6873
$ws.Range("L134").Value = 280679.46
$ws.Range("M134").Value = -1896.6873
$ws.Range("N134").Value = -285749.46

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2818.0454
$ws.Range("I31").Value = 2581.0938
$ws.Range("J31").Value = 3449.9167
$ws.Range("K31").Value = 2581.0938
$ws.Range("L31").Value = 3449.9167
$ws.Range("M31").Value = -2286.0938
$ws.Range("N31").Value = -4039.9167
$ws.Range("H34").Value = 2818.0454
$ws.Range("I34").Value = 2581.0938
$ws.Range("J34").Value = 3449.9167
$ws.Range("K34").Value = 2581.0938
$ws.Range("L34").Value = 3449.9167
$ws.Range("M34").Value = -2379.0938
$ws.Range("N34").Value = -3853.9167
$ws.Range("H99").Value = 2643.2
$ws.Range("I99").Value = 2241.0527
$ws.Range("J99").Value = 3916.6667
$ws.Range("K99").Value = 2241.0527
$ws.Range("L99").Value = 3916.6667
$ws.Range("M99").Value = -743.0527000000002
$ws.Range("N99").Value = -6912.6667
$ws.Range("H126").Value = 2643.2
$ws.Range("I126").Value = 2241.0527
$ws.Range("J126").Value = 3916.6667
$ws.Range("K126").Value = 6723.158100000001
$ws.Range("L126").Value = 11750.0001
$ws.Range("M126").Value = -4253.158100000001
$ws.Range("N126").Value = -16690.0001
$ws.Range("H132").Value = 1647.4626
$ws.Range("I132").Value = 933.1892
$ws.Range("J132").Value = 2528.4
$ws.Range("K132").Value = 2799.5676
$ws.Range("L132").Value = 7585.200000000001
$ws.Range("M132").Value = -269.5676000000003
$ws.Range("N132").Value = -12645.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2960.5
$ws.Range("I56").Value = 2960.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 2960.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2430.5
$ws.Range("H122").Value = 568.18604
$ws.Range("I122").Value = 258.57144
$ws.Range("J122").Value = 717.65515
$ws.Range("K122").Value = 2327.14296
$ws.Range("L122").Value = 6458.896350000001
$ws.Range("M122").Value = 122.8570399999999
$ws.Range("N122").Value = -11358.89635
$ws.Range("H129").Value = 3055.037
$ws.Range("I129").Value = 1046.6666
$ws.Range("J129").Value = 4059.2222
$ws.Range("K129").Value = 3139.9998
$ws.Range("L129").Value = 12177.6666
$ws.Range("M129").Value = 1860.0002
$ws.Range("N129").Value = -22177.6666
$ws.Range("H131").Value = 20329.885
$ws.Range("I131").Value = 202292
$ws.Range("J131").Value = 972.21277
$ws.Range("K131").Value = 606876
$ws.Range("L131").Value = 2916.63831
$ws.Range("M131").Value = -601836
$ws.Range("N131").Value = -12996.63831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2473.838
$ws.Range("I132").Value = 2250.2917
$ws.Range("J132").Value = 2886.5386
$ws.Range("K132").Value = 6750.875100000001
$ws.Range("L132").Value = 8659.6158
$ws.Range("M132").Value = -4220.875100000001
$ws.Range("N132").Value = -13719.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2210.1333
$ws.Range("I68").Value = 2507.4
$ws.Range("J68").Value = 1615.6
$ws.Range("K68").Value = 2507.4
$ws.Range("L68").Value = 1615.6
$ws.Range("M68").Value = -1758.4
$ws.Range("N68").Value = -3113.6
$ws.Range("H71").Value = 2210.1333
$ws.Range("I71").Value = 2507.4
$ws.Range("J71").Value = 1615.6
$ws.Range("K71").Value = 12537
$ws.Range("L71").Value = 8078
$ws.Range("M71").Value = -8793
$ws.Range("N71").Value = -15566
$ws.Range("H132").Value = 2035.9149
$ws.Range("I132").Value = 1754.5
$ws.Range("J132").Value = 4399.8
$ws.Range("K132").Value = 5263.5
$ws.Range("L132").Value = 13199.4
$ws.Range("M132").Value = -2733.5
$ws.Range("N132").Value = -18259.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1053.875
$ws.Range("I107").Value = 1381.2
$ws.Range("J107").Value = 508.33334
$ws.Range("K107").Value = 4143.6
$ws.Range("L107").Value = 1525.00002
$ws.Range("M107").Value = -2223.6
$ws.Range("N107").Value = -5365.000019999999
$ws.Range("H132").Value = 641.1774
$ws.Range("I132").Value = 532.7091
$ws.Range("J132").Value = 1493.4286
$ws.Range("K132").Value = 1598.1273
$ws.Range("L132").Value = 4480.2858
$ws.Range("M132").Value = 931.8726999999999
$ws.Range("N132").Value = -9540.2858
$ws.Range("H136").Value = 1056.5151
$ws.Range("I136").Value = 1297.85
$ws.Range("J136").Value = 685.2308
$ws.Range("K136").Value = 3893.55
$ws.Range("L136").Value = 2055.6924
$ws.Range("M136").Value = -1343.55
$ws.Range("N136").Value = -7155.6924
